# Auto-generated Excel COM-interop script
# Applies the "Horarios actualizados Linea 141 - 271" update:
# refreshes the scrape timestamp/row-count headers and merges newly
# scraped rows (timestamp 12:01:50) into the sorted schedule tables
# on all three worksheets (LP1912, LP1912-215, 6203-6173).

$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item('LP1912')
$ws.Cells.Item(2,1).Value = 'Última actualización: 12:01:50'
$ws.Cells.Item(3,1).Value = 'Total filas: 223'

$ws.Cells.Item(47,1).Value = '05:49:40'; $ws.Cells.Item(47,2).Value = '07:04'; $ws.Cells.Item(47,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(47,4).Value = 75; $ws.Cells.Item(47,5).Value = 'LP1912'
$ws.Cells.Item(48,1).Value = '05:18:56'; $ws.Cells.Item(48,2).Value = '07:04'; $ws.Cells.Item(48,3).Value = '15_ABASTO'; $ws.Cells.Item(48,4).Value = 106; $ws.Cells.Item(48,5).Value = 'LP1912'
$ws.Cells.Item(49,1).Value = '05:49:40'; $ws.Cells.Item(49,2).Value = '07:05'; $ws.Cells.Item(49,3).Value = '15_ABASTO'; $ws.Cells.Item(49,4).Value = 76; $ws.Cells.Item(49,5).Value = 'LP1912'
$ws.Cells.Item(50,1).Value = '05:18:56'; $ws.Cells.Item(50,2).Value = '07:06'; $ws.Cells.Item(50,3).Value = '225_GOMEZ'; $ws.Cells.Item(50,4).Value = 108; $ws.Cells.Item(50,5).Value = 'LP1912'
$ws.Cells.Item(51,1).Value = '05:49:40'; $ws.Cells.Item(51,2).Value = '07:07'; $ws.Cells.Item(51,3).Value = '225_GOMEZ'; $ws.Cells.Item(51,4).Value = 78; $ws.Cells.Item(51,5).Value = 'LP1912'
$ws.Cells.Item(52,1).Value = '05:18:56'; $ws.Cells.Item(52,2).Value = '07:11'; $ws.Cells.Item(52,3).Value = '215A_EL PATO'; $ws.Cells.Item(52,4).Value = 113; $ws.Cells.Item(52,5).Value = 'LP1912'
$ws.Cells.Item(53,1).Value = '06:15:04'; $ws.Cells.Item(53,2).Value = '07:12'; $ws.Cells.Item(53,3).Value = '215A_EL PATO'; $ws.Cells.Item(53,4).Value = 57; $ws.Cells.Item(53,5).Value = 'LP1912'
$ws.Cells.Item(54,1).Value = '05:18:56'; $ws.Cells.Item(54,2).Value = '07:15'; $ws.Cells.Item(54,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(54,4).Value = 117; $ws.Cells.Item(54,5).Value = 'LP1912'
$ws.Cells.Item(55,1).Value = '06:15:04'; $ws.Cells.Item(55,2).Value = '07:16'; $ws.Cells.Item(55,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(55,4).Value = 61; $ws.Cells.Item(55,5).Value = 'LP1912'
$ws.Cells.Item(56,1).Value = '06:43:40'; $ws.Cells.Item(56,2).Value = '07:16'; $ws.Cells.Item(56,3).Value = '16_SANTA ANA'; $ws.Cells.Item(56,4).Value = 33; $ws.Cells.Item(56,5).Value = 'LP1912'
$ws.Cells.Item(57,1).Value = '07:20:40'; $ws.Cells.Item(57,2).Value = '07:20'; $ws.Cells.Item(57,3).Value = '10_OLMOS'; $ws.Cells.Item(57,4).Value = 0; $ws.Cells.Item(57,5).Value = 'LP1912'
$ws.Cells.Item(58,1).Value = '05:49:40'; $ws.Cells.Item(58,2).Value = '07:21'; $ws.Cells.Item(58,3).Value = '26_HERNANDEZ'; $ws.Cells.Item(58,4).Value = 92; $ws.Cells.Item(58,5).Value = 'LP1912'
$ws.Cells.Item(59,1).Value = '06:15:04'; $ws.Cells.Item(59,2).Value = '07:23'; $ws.Cells.Item(59,3).Value = '10_OLMOS'; $ws.Cells.Item(59,4).Value = 68; $ws.Cells.Item(59,5).Value = 'LP1912'
$ws.Cells.Item(60,1).Value = '05:49:40'; $ws.Cells.Item(60,2).Value = '07:29'; $ws.Cells.Item(60,3).Value = '10_OLMOS'; $ws.Cells.Item(60,4).Value = 100; $ws.Cells.Item(60,5).Value = 'LP1912'
$ws.Cells.Item(61,1).Value = '05:49:40'; $ws.Cells.Item(61,2).Value = '07:31'; $ws.Cells.Item(61,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(61,4).Value = 102; $ws.Cells.Item(61,5).Value = 'LP1912'
$ws.Cells.Item(62,1).Value = '05:49:40'; $ws.Cells.Item(62,2).Value = '07:32'; $ws.Cells.Item(62,3).Value = '84_COLONIA URQUIZA-ESC 49'; $ws.Cells.Item(62,4).Value = 103; $ws.Cells.Item(62,5).Value = 'LP1912'
$ws.Cells.Item(63,1).Value = '06:15:04'; $ws.Cells.Item(63,2).Value = '07:32'; $ws.Cells.Item(63,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(63,4).Value = 77; $ws.Cells.Item(63,5).Value = 'LP1912'
$ws.Cells.Item(64,1).Value = '07:20:40'; $ws.Cells.Item(64,2).Value = '07:34'; $ws.Cells.Item(64,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(64,4).Value = 14; $ws.Cells.Item(64,5).Value = 'LP1912'
$ws.Cells.Item(65,1).Value = '05:49:40'; $ws.Cells.Item(65,2).Value = '07:36'; $ws.Cells.Item(65,3).Value = '27_EL RETIRO'; $ws.Cells.Item(65,4).Value = 107; $ws.Cells.Item(65,5).Value = 'LP1912'
$ws.Cells.Item(66,1).Value = '06:15:04'; $ws.Cells.Item(66,2).Value = '07:37'; $ws.Cells.Item(66,3).Value = '27_EL RETIRO'; $ws.Cells.Item(66,4).Value = 82; $ws.Cells.Item(66,5).Value = 'LP1912'
$ws.Cells.Item(67,1).Value = '05:49:40'; $ws.Cells.Item(67,2).Value = '07:39'; $ws.Cells.Item(67,3).Value = '10_OLMOS'; $ws.Cells.Item(67,4).Value = 110; $ws.Cells.Item(67,5).Value = 'LP1912'
$ws.Cells.Item(68,1).Value = '07:20:40'; $ws.Cells.Item(68,2).Value = '07:46'; $ws.Cells.Item(68,3).Value = '16_SANTA ANA'; $ws.Cells.Item(68,4).Value = 26; $ws.Cells.Item(68,5).Value = 'LP1912'
$ws.Cells.Item(69,1).Value = '06:43:40'; $ws.Cells.Item(69,2).Value = '07:47'; $ws.Cells.Item(69,3).Value = '14_ABASTO'; $ws.Cells.Item(69,4).Value = 64; $ws.Cells.Item(69,5).Value = 'LP1912'
$ws.Cells.Item(70,1).Value = '06:15:04'; $ws.Cells.Item(70,2).Value = '07:48'; $ws.Cells.Item(70,3).Value = '14_ABASTO'; $ws.Cells.Item(70,4).Value = 93; $ws.Cells.Item(70,5).Value = 'LP1912'
$ws.Cells.Item(71,1).Value = '06:43:40'; $ws.Cells.Item(71,2).Value = '07:51'; $ws.Cells.Item(71,3).Value = '215D_EL PATO'; $ws.Cells.Item(71,4).Value = 68; $ws.Cells.Item(71,5).Value = 'LP1912'
$ws.Cells.Item(72,1).Value = '06:15:04'; $ws.Cells.Item(72,2).Value = '07:52'; $ws.Cells.Item(72,3).Value = '215D_EL PATO'; $ws.Cells.Item(72,4).Value = 97; $ws.Cells.Item(72,5).Value = 'LP1912'
$ws.Cells.Item(73,1).Value = '07:47:32'; $ws.Cells.Item(73,2).Value = '07:55'; $ws.Cells.Item(73,3).Value = '10_OLMOS'; $ws.Cells.Item(73,4).Value = 8; $ws.Cells.Item(73,5).Value = 'LP1912'
$ws.Cells.Item(74,1).Value = '07:20:40'; $ws.Cells.Item(74,2).Value = '07:58'; $ws.Cells.Item(74,3).Value = '16_SANTA ANA'; $ws.Cells.Item(74,4).Value = 38; $ws.Cells.Item(74,5).Value = 'LP1912'
$ws.Cells.Item(75,1).Value = '07:20:40'; $ws.Cells.Item(75,2).Value = '07:59'; $ws.Cells.Item(75,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(75,4).Value = 39; $ws.Cells.Item(75,5).Value = 'LP1912'
$ws.Cells.Item(76,1).Value = '06:15:04'; $ws.Cells.Item(76,2).Value = '08:01'; $ws.Cells.Item(76,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(76,4).Value = 106; $ws.Cells.Item(76,5).Value = 'LP1912'
$ws.Cells.Item(77,1).Value = '06:43:40'; $ws.Cells.Item(77,2).Value = '08:03'; $ws.Cells.Item(77,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(77,4).Value = 80; $ws.Cells.Item(77,5).Value = 'LP1912'
$ws.Cells.Item(78,1).Value = '07:20:40'; $ws.Cells.Item(78,2).Value = '08:03'; $ws.Cells.Item(78,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(78,4).Value = 43; $ws.Cells.Item(78,5).Value = 'LP1912'
$ws.Cells.Item(79,1).Value = '06:57:30'; $ws.Cells.Item(79,2).Value = '08:06'; $ws.Cells.Item(79,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(79,4).Value = 69; $ws.Cells.Item(79,5).Value = 'LP1912'
$ws.Cells.Item(80,1).Value = '07:47:32'; $ws.Cells.Item(80,2).Value = '08:10'; $ws.Cells.Item(80,3).Value = '16_SANTA ANA'; $ws.Cells.Item(80,4).Value = 23; $ws.Cells.Item(80,5).Value = 'LP1912'
$ws.Cells.Item(81,1).Value = '06:15:04'; $ws.Cells.Item(81,2).Value = '08:12'; $ws.Cells.Item(81,3).Value = '15_ABASTO'; $ws.Cells.Item(81,4).Value = 117; $ws.Cells.Item(81,5).Value = 'LP1912'
$ws.Cells.Item(82,1).Value = '07:47:32'; $ws.Cells.Item(82,2).Value = '08:13'; $ws.Cells.Item(82,3).Value = '10_OLMOS'; $ws.Cells.Item(82,4).Value = 26; $ws.Cells.Item(82,5).Value = 'LP1912'
$ws.Cells.Item(83,1).Value = '07:47:32'; $ws.Cells.Item(83,2).Value = '08:16'; $ws.Cells.Item(83,3).Value = '26_HERNANDEZ'; $ws.Cells.Item(83,4).Value = 29; $ws.Cells.Item(83,5).Value = 'LP1912'
$ws.Cells.Item(84,1).Value = '06:43:40'; $ws.Cells.Item(84,2).Value = '08:21'; $ws.Cells.Item(84,3).Value = '26_HERNANDEZ'; $ws.Cells.Item(84,4).Value = 98; $ws.Cells.Item(84,5).Value = 'LP1912'
$ws.Cells.Item(85,1).Value = '06:43:40'; $ws.Cells.Item(85,2).Value = '08:22'; $ws.Cells.Item(85,3).Value = '16_P MOR-SANTA ANA'; $ws.Cells.Item(85,4).Value = 99; $ws.Cells.Item(85,5).Value = 'LP1912'
$ws.Cells.Item(86,1).Value = '06:43:40'; $ws.Cells.Item(86,2).Value = '08:23'; $ws.Cells.Item(86,3).Value = '215B_EL PATO'; $ws.Cells.Item(86,4).Value = 100; $ws.Cells.Item(86,5).Value = 'LP1912'
$ws.Cells.Item(87,1).Value = '06:43:40'; $ws.Cells.Item(87,2).Value = '08:27'; $ws.Cells.Item(87,3).Value = '84_COLONIA URQUIZA-ESC 49'; $ws.Cells.Item(87,4).Value = 104; $ws.Cells.Item(87,5).Value = 'LP1912'
$ws.Cells.Item(88,1).Value = '07:47:32'; $ws.Cells.Item(88,2).Value = '08:31'; $ws.Cells.Item(88,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(88,4).Value = 44; $ws.Cells.Item(88,5).Value = 'LP1912'
$ws.Cells.Item(89,1).Value = '07:59:28'; $ws.Cells.Item(89,2).Value = '08:33'; $ws.Cells.Item(89,3).Value = '10_OLMOS'; $ws.Cells.Item(89,4).Value = 34; $ws.Cells.Item(89,5).Value = 'LP1912'
$ws.Cells.Item(90,1).Value = '07:59:28'; $ws.Cells.Item(90,2).Value = '08:34'; $ws.Cells.Item(90,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(90,4).Value = 35; $ws.Cells.Item(90,5).Value = 'LP1912'
$ws.Cells.Item(91,1).Value = '07:59:28'; $ws.Cells.Item(91,2).Value = '08:39'; $ws.Cells.Item(91,3).Value = '26_HERNANDEZ'; $ws.Cells.Item(91,4).Value = 40; $ws.Cells.Item(91,5).Value = 'LP1912'
$ws.Cells.Item(92,1).Value = '06:43:40'; $ws.Cells.Item(92,2).Value = '08:42'; $ws.Cells.Item(92,3).Value = '81_EL PELIGRO'; $ws.Cells.Item(92,4).Value = 119; $ws.Cells.Item(92,5).Value = 'LP1912'
$ws.Cells.Item(93,1).Value = '07:20:40'; $ws.Cells.Item(93,2).Value = '08:43'; $ws.Cells.Item(93,3).Value = '14_ABASTO'; $ws.Cells.Item(93,4).Value = 83; $ws.Cells.Item(93,5).Value = 'LP1912'
$ws.Cells.Item(94,1).Value = '06:57:30'; $ws.Cells.Item(94,2).Value = '08:54'; $ws.Cells.Item(94,3).Value = '17_ROMERO'; $ws.Cells.Item(94,4).Value = 117; $ws.Cells.Item(94,5).Value = 'LP1912'
$ws.Cells.Item(95,1).Value = '08:57:13'; $ws.Cells.Item(95,2).Value = '08:59'; $ws.Cells.Item(95,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(95,4).Value = 2; $ws.Cells.Item(95,5).Value = 'LP1912'
$ws.Cells.Item(96,1).Value = '08:21:50'; $ws.Cells.Item(96,2).Value = '09:01'; $ws.Cells.Item(96,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(96,4).Value = 40; $ws.Cells.Item(96,5).Value = 'LP1912'
$ws.Cells.Item(97,1).Value = '07:20:40'; $ws.Cells.Item(97,2).Value = '09:01'; $ws.Cells.Item(97,3).Value = '215A_EL PATO'; $ws.Cells.Item(97,4).Value = 101; $ws.Cells.Item(97,5).Value = 'LP1912'
$ws.Cells.Item(98,1).Value = '08:57:13'; $ws.Cells.Item(98,2).Value = '09:02'; $ws.Cells.Item(98,3).Value = '215A_EL PATO'; $ws.Cells.Item(98,4).Value = 5; $ws.Cells.Item(98,5).Value = 'LP1912'
$ws.Cells.Item(99,1).Value = '07:59:28'; $ws.Cells.Item(99,2).Value = '09:03'; $ws.Cells.Item(99,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(99,4).Value = 64; $ws.Cells.Item(99,5).Value = 'LP1912'
$ws.Cells.Item(100,1).Value = '08:39:44'; $ws.Cells.Item(100,2).Value = '09:04'; $ws.Cells.Item(100,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(100,4).Value = 25; $ws.Cells.Item(100,5).Value = 'LP1912'
$ws.Cells.Item(101,1).Value = '08:57:13'; $ws.Cells.Item(101,2).Value = '09:05'; $ws.Cells.Item(101,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(101,4).Value = 8; $ws.Cells.Item(101,5).Value = 'LP1912'
$ws.Cells.Item(102,1).Value = '08:21:50'; $ws.Cells.Item(102,2).Value = '09:07'; $ws.Cells.Item(102,3).Value = '26_HERNANDEZ'; $ws.Cells.Item(102,4).Value = 46; $ws.Cells.Item(102,5).Value = 'LP1912'
$ws.Cells.Item(103,1).Value = '07:20:40'; $ws.Cells.Item(103,2).Value = '09:10'; $ws.Cells.Item(103,3).Value = '16_P MOR-SANTA ANA'; $ws.Cells.Item(103,4).Value = 110; $ws.Cells.Item(103,5).Value = 'LP1912'
$ws.Cells.Item(104,1).Value = '08:57:13'; $ws.Cells.Item(104,2).Value = '09:11'; $ws.Cells.Item(104,3).Value = '16_P MOR-SANTA ANA'; $ws.Cells.Item(104,4).Value = 14; $ws.Cells.Item(104,5).Value = 'LP1912'
$ws.Cells.Item(105,1).Value = '08:21:50'; $ws.Cells.Item(105,2).Value = '09:13'; $ws.Cells.Item(105,3).Value = '10_OLMOS'; $ws.Cells.Item(105,4).Value = 52; $ws.Cells.Item(105,5).Value = 'LP1912'
$ws.Cells.Item(106,1).Value = '07:20:40'; $ws.Cells.Item(106,2).Value = '09:16'; $ws.Cells.Item(106,3).Value = '27_EL RETIRO'; $ws.Cells.Item(106,4).Value = 116; $ws.Cells.Item(106,5).Value = 'LP1912'
$ws.Cells.Item(107,1).Value = '08:57:13'; $ws.Cells.Item(107,2).Value = '09:17'; $ws.Cells.Item(107,3).Value = '27_EL RETIRO'; $ws.Cells.Item(107,4).Value = 20; $ws.Cells.Item(107,5).Value = 'LP1912'
$ws.Cells.Item(108,1).Value = '08:21:50'; $ws.Cells.Item(108,2).Value = '09:21'; $ws.Cells.Item(108,3).Value = '26_HERNANDEZ'; $ws.Cells.Item(108,4).Value = 60; $ws.Cells.Item(108,5).Value = 'LP1912'
$ws.Cells.Item(109,1).Value = '07:59:28'; $ws.Cells.Item(109,2).Value = '09:22'; $ws.Cells.Item(109,3).Value = '16_SANTA ANA'; $ws.Cells.Item(109,4).Value = 83; $ws.Cells.Item(109,5).Value = 'LP1912'
$ws.Cells.Item(110,1).Value = '07:47:32'; $ws.Cells.Item(110,2).Value = '09:22'; $ws.Cells.Item(110,3).Value = '17_ROMERO'; $ws.Cells.Item(110,4).Value = 95; $ws.Cells.Item(110,5).Value = 'LP1912'
$ws.Cells.Item(111,1).Value = '08:57:13'; $ws.Cells.Item(111,2).Value = '09:23'; $ws.Cells.Item(111,3).Value = '16_SANTA ANA'; $ws.Cells.Item(111,4).Value = 26; $ws.Cells.Item(111,5).Value = 'LP1912'
$ws.Cells.Item(112,1).Value = '07:47:32'; $ws.Cells.Item(112,2).Value = '09:23'; $ws.Cells.Item(112,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(112,4).Value = 96; $ws.Cells.Item(112,5).Value = 'LP1912'
$ws.Cells.Item(113,1).Value = '08:57:13'; $ws.Cells.Item(113,2).Value = '09:24'; $ws.Cells.Item(113,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(113,4).Value = 27; $ws.Cells.Item(113,5).Value = 'LP1912'
$ws.Cells.Item(114,1).Value = '08:21:50'; $ws.Cells.Item(114,2).Value = '09:29'; $ws.Cells.Item(114,3).Value = '16_SANTA ANA'; $ws.Cells.Item(114,4).Value = 68; $ws.Cells.Item(114,5).Value = 'LP1912'
$ws.Cells.Item(115,1).Value = '07:47:32'; $ws.Cells.Item(115,2).Value = '09:32'; $ws.Cells.Item(115,3).Value = '15_ABASTO'; $ws.Cells.Item(115,4).Value = 105; $ws.Cells.Item(115,5).Value = 'LP1912'
$ws.Cells.Item(116,1).Value = '07:47:32'; $ws.Cells.Item(116,2).Value = '09:33'; $ws.Cells.Item(116,3).Value = '10_OLMOS'; $ws.Cells.Item(116,4).Value = 106; $ws.Cells.Item(116,5).Value = 'LP1912'
$ws.Cells.Item(117,1).Value = '08:39:44'; $ws.Cells.Item(117,2).Value = '09:34'; $ws.Cells.Item(117,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(117,4).Value = 55; $ws.Cells.Item(117,5).Value = 'LP1912'
$ws.Cells.Item(118,1).Value = '08:39:44'; $ws.Cells.Item(118,2).Value = '09:34'; $ws.Cells.Item(118,3).Value = '16_SANTA ANA'; $ws.Cells.Item(118,4).Value = 55; $ws.Cells.Item(118,5).Value = 'LP1912'
$ws.Cells.Item(119,1).Value = '08:57:13'; $ws.Cells.Item(119,2).Value = '09:35'; $ws.Cells.Item(119,3).Value = '16_SANTA ANA'; $ws.Cells.Item(119,4).Value = 38; $ws.Cells.Item(119,5).Value = 'LP1912'
$ws.Cells.Item(120,1).Value = '08:57:13'; $ws.Cells.Item(120,2).Value = '09:35'; $ws.Cells.Item(120,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(120,4).Value = 38; $ws.Cells.Item(120,5).Value = 'LP1912'
$ws.Cells.Item(121,1).Value = '09:38:09'; $ws.Cells.Item(121,2).Value = '09:41'; $ws.Cells.Item(121,3).Value = '14_ABASTO'; $ws.Cells.Item(121,4).Value = 3; $ws.Cells.Item(121,5).Value = 'LP1912'
$ws.Cells.Item(122,1).Value = '08:21:50'; $ws.Cells.Item(122,2).Value = '09:41'; $ws.Cells.Item(122,3).Value = '215C_EL PATO'; $ws.Cells.Item(122,4).Value = 80; $ws.Cells.Item(122,5).Value = 'LP1912'
$ws.Cells.Item(123,1).Value = '09:38:09'; $ws.Cells.Item(123,2).Value = '09:41'; $ws.Cells.Item(123,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(123,4).Value = 3; $ws.Cells.Item(123,5).Value = 'LP1912'
$ws.Cells.Item(124,1).Value = '07:47:32'; $ws.Cells.Item(124,2).Value = '09:42'; $ws.Cells.Item(124,3).Value = '215C_EL PATO'; $ws.Cells.Item(124,4).Value = 115; $ws.Cells.Item(124,5).Value = 'LP1912'
$ws.Cells.Item(125,1).Value = '07:47:32'; $ws.Cells.Item(125,2).Value = '09:43'; $ws.Cells.Item(125,3).Value = '14_ABASTO'; $ws.Cells.Item(125,4).Value = 116; $ws.Cells.Item(125,5).Value = 'LP1912'
$ws.Cells.Item(126,1).Value = '08:57:13'; $ws.Cells.Item(126,2).Value = '09:44'; $ws.Cells.Item(126,3).Value = '14_ABASTO'; $ws.Cells.Item(126,4).Value = 47; $ws.Cells.Item(126,5).Value = 'LP1912'
$ws.Cells.Item(127,1).Value = '09:38:09'; $ws.Cells.Item(127,2).Value = '09:47'; $ws.Cells.Item(127,3).Value = '16_SANTA ANA'; $ws.Cells.Item(127,4).Value = 9; $ws.Cells.Item(127,5).Value = 'LP1912'
$ws.Cells.Item(128,1).Value = '08:49:51'; $ws.Cells.Item(128,2).Value = '09:52'; $ws.Cells.Item(128,3).Value = '15_ABASTO'; $ws.Cells.Item(128,4).Value = 63; $ws.Cells.Item(128,5).Value = 'LP1912'
$ws.Cells.Item(129,1).Value = '08:49:51'; $ws.Cells.Item(129,2).Value = '09:53'; $ws.Cells.Item(129,3).Value = '10_OLMOS'; $ws.Cells.Item(129,4).Value = 64; $ws.Cells.Item(129,5).Value = 'LP1912'
$ws.Cells.Item(130,1).Value = '09:38:09'; $ws.Cells.Item(130,2).Value = '09:59'; $ws.Cells.Item(130,3).Value = '16_SANTA ANA'; $ws.Cells.Item(130,4).Value = 21; $ws.Cells.Item(130,5).Value = 'LP1912'
$ws.Cells.Item(131,1).Value = '09:38:09'; $ws.Cells.Item(131,2).Value = '10:04'; $ws.Cells.Item(131,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(131,4).Value = 26; $ws.Cells.Item(131,5).Value = 'LP1912'
$ws.Cells.Item(132,1).Value = '09:38:09'; $ws.Cells.Item(132,2).Value = '10:05'; $ws.Cells.Item(132,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(132,4).Value = 27; $ws.Cells.Item(132,5).Value = 'LP1912'
$ws.Cells.Item(133,1).Value = '08:39:44'; $ws.Cells.Item(133,2).Value = '10:06'; $ws.Cells.Item(133,3).Value = '10_OLMOS'; $ws.Cells.Item(133,4).Value = 87; $ws.Cells.Item(133,5).Value = 'LP1912'
$ws.Cells.Item(134,1).Value = '08:21:50'; $ws.Cells.Item(134,2).Value = '10:10'; $ws.Cells.Item(134,3).Value = '16_P MOR-SANTA ANA'; $ws.Cells.Item(134,4).Value = 109; $ws.Cells.Item(134,5).Value = 'LP1912'
$ws.Cells.Item(135,1).Value = '08:57:13'; $ws.Cells.Item(135,2).Value = '10:11'; $ws.Cells.Item(135,3).Value = '16_P MOR-SANTA ANA'; $ws.Cells.Item(135,4).Value = 74; $ws.Cells.Item(135,5).Value = 'LP1912'
$ws.Cells.Item(136,1).Value = '08:21:50'; $ws.Cells.Item(136,2).Value = '10:12'; $ws.Cells.Item(136,3).Value = '15_ABASTO'; $ws.Cells.Item(136,4).Value = 111; $ws.Cells.Item(136,5).Value = 'LP1912'
$ws.Cells.Item(137,1).Value = '09:38:09'; $ws.Cells.Item(137,2).Value = '10:13'; $ws.Cells.Item(137,3).Value = '10_OLMOS'; $ws.Cells.Item(137,4).Value = 35; $ws.Cells.Item(137,5).Value = 'LP1912'
$ws.Cells.Item(138,1).Value = '08:49:51'; $ws.Cells.Item(138,2).Value = '10:20'; $ws.Cells.Item(138,3).Value = '26_HERNANDEZ'; $ws.Cells.Item(138,4).Value = 91; $ws.Cells.Item(138,5).Value = 'LP1912'
$ws.Cells.Item(139,1).Value = '08:39:44'; $ws.Cells.Item(139,2).Value = '10:21'; $ws.Cells.Item(139,3).Value = '26_HERNANDEZ'; $ws.Cells.Item(139,4).Value = 102; $ws.Cells.Item(139,5).Value = 'LP1912'
$ws.Cells.Item(140,1).Value = '08:39:44'; $ws.Cells.Item(140,2).Value = '10:22'; $ws.Cells.Item(140,3).Value = '17_ROMERO'; $ws.Cells.Item(140,4).Value = 103; $ws.Cells.Item(140,5).Value = 'LP1912'
$ws.Cells.Item(141,1).Value = '09:38:09'; $ws.Cells.Item(141,2).Value = '10:24'; $ws.Cells.Item(141,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(141,4).Value = 46; $ws.Cells.Item(141,5).Value = 'LP1912'
$ws.Cells.Item(142,1).Value = '08:39:44'; $ws.Cells.Item(142,2).Value = '10:26'; $ws.Cells.Item(142,3).Value = '215A_EL PATO'; $ws.Cells.Item(142,4).Value = 107; $ws.Cells.Item(142,5).Value = 'LP1912'
$ws.Cells.Item(143,1).Value = '08:57:13'; $ws.Cells.Item(143,2).Value = '10:27'; $ws.Cells.Item(143,3).Value = '215A_EL PATO'; $ws.Cells.Item(143,4).Value = 90; $ws.Cells.Item(143,5).Value = 'LP1912'
$ws.Cells.Item(144,1).Value = '10:26:41'; $ws.Cells.Item(144,2).Value = '10:33'; $ws.Cells.Item(144,3).Value = '10_OLMOS'; $ws.Cells.Item(144,4).Value = 7; $ws.Cells.Item(144,5).Value = 'LP1912'
$ws.Cells.Item(145,1).Value = '10:26:41'; $ws.Cells.Item(145,2).Value = '10:34'; $ws.Cells.Item(145,3).Value = '16_SANTA ANA'; $ws.Cells.Item(145,4).Value = 8; $ws.Cells.Item(145,5).Value = 'LP1912'
$ws.Cells.Item(146,1).Value = '10:26:41'; $ws.Cells.Item(146,2).Value = '10:34'; $ws.Cells.Item(146,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(146,4).Value = 8; $ws.Cells.Item(146,5).Value = 'LP1912'
$ws.Cells.Item(147,1).Value = '08:49:51'; $ws.Cells.Item(147,2).Value = '10:41'; $ws.Cells.Item(147,3).Value = '17_ROMERO'; $ws.Cells.Item(147,4).Value = 112; $ws.Cells.Item(147,5).Value = 'LP1912'
$ws.Cells.Item(148,1).Value = '08:57:13'; $ws.Cells.Item(148,2).Value = '10:42'; $ws.Cells.Item(148,3).Value = '17_ROMERO'; $ws.Cells.Item(148,4).Value = 105; $ws.Cells.Item(148,5).Value = 'LP1912'
$ws.Cells.Item(149,1).Value = '08:49:51'; $ws.Cells.Item(149,2).Value = '10:43'; $ws.Cells.Item(149,3).Value = '14_ABASTO'; $ws.Cells.Item(149,4).Value = 114; $ws.Cells.Item(149,5).Value = 'LP1912'
$ws.Cells.Item(150,1).Value = '08:57:13'; $ws.Cells.Item(150,2).Value = '10:44'; $ws.Cells.Item(150,3).Value = '14_ABASTO'; $ws.Cells.Item(150,4).Value = 107; $ws.Cells.Item(150,5).Value = 'LP1912'
$ws.Cells.Item(151,1).Value = '10:26:41'; $ws.Cells.Item(151,2).Value = '10:46'; $ws.Cells.Item(151,3).Value = '16_SANTA ANA'; $ws.Cells.Item(151,4).Value = 20; $ws.Cells.Item(151,5).Value = 'LP1912'
$ws.Cells.Item(152,1).Value = '10:26:41'; $ws.Cells.Item(152,2).Value = '10:52'; $ws.Cells.Item(152,3).Value = '15_ABASTO'; $ws.Cells.Item(152,4).Value = 26; $ws.Cells.Item(152,5).Value = 'LP1912'
$ws.Cells.Item(153,1).Value = '10:26:41'; $ws.Cells.Item(153,2).Value = '10:53'; $ws.Cells.Item(153,3).Value = '10_OLMOS'; $ws.Cells.Item(153,4).Value = 27; $ws.Cells.Item(153,5).Value = 'LP1912'
$ws.Cells.Item(154,1).Value = '10:26:41'; $ws.Cells.Item(154,2).Value = '10:56'; $ws.Cells.Item(154,3).Value = '27_EL RETIRO'; $ws.Cells.Item(154,4).Value = 30; $ws.Cells.Item(154,5).Value = 'LP1912'
$ws.Cells.Item(155,1).Value = '10:56:30'; $ws.Cells.Item(155,2).Value = '10:57'; $ws.Cells.Item(155,3).Value = '16_SANTA ANA'; $ws.Cells.Item(155,4).Value = 1; $ws.Cells.Item(155,5).Value = 'LP1912'
$ws.Cells.Item(156,1).Value = '09:38:09'; $ws.Cells.Item(156,2).Value = '10:58'; $ws.Cells.Item(156,3).Value = '27_EL RETIRO'; $ws.Cells.Item(156,4).Value = 80; $ws.Cells.Item(156,5).Value = 'LP1912'
$ws.Cells.Item(157,1).Value = '10:26:41'; $ws.Cells.Item(157,2).Value = '11:01'; $ws.Cells.Item(157,3).Value = '215C_EL PATO'; $ws.Cells.Item(157,4).Value = 35; $ws.Cells.Item(157,5).Value = 'LP1912'
$ws.Cells.Item(158,1).Value = '09:38:09'; $ws.Cells.Item(158,2).Value = '11:02'; $ws.Cells.Item(158,3).Value = '215C_EL PATO'; $ws.Cells.Item(158,4).Value = 84; $ws.Cells.Item(158,5).Value = 'LP1912'
$ws.Cells.Item(159,1).Value = '10:26:41'; $ws.Cells.Item(159,2).Value = '11:03'; $ws.Cells.Item(159,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(159,4).Value = 37; $ws.Cells.Item(159,5).Value = 'LP1912'
$ws.Cells.Item(160,1).Value = '10:26:41'; $ws.Cells.Item(160,2).Value = '11:04'; $ws.Cells.Item(160,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(160,4).Value = 38; $ws.Cells.Item(160,5).Value = 'LP1912'
$ws.Cells.Item(161,1).Value = '10:26:41'; $ws.Cells.Item(161,2).Value = '11:06'; $ws.Cells.Item(161,3).Value = '16_P MOR-167 Y 521'; $ws.Cells.Item(161,4).Value = 40; $ws.Cells.Item(161,5).Value = 'LP1912'
$ws.Cells.Item(162,1).Value = '09:38:09'; $ws.Cells.Item(162,2).Value = '11:07'; $ws.Cells.Item(162,3).Value = '16_P MOR-167 Y 521'; $ws.Cells.Item(162,4).Value = 89; $ws.Cells.Item(162,5).Value = 'LP1912'
$ws.Cells.Item(163,1).Value = '10:56:30'; $ws.Cells.Item(163,2).Value = '11:08'; $ws.Cells.Item(163,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(163,4).Value = 12; $ws.Cells.Item(163,5).Value = 'LP1912'
$ws.Cells.Item(164,1).Value = '10:26:41'; $ws.Cells.Item(164,2).Value = '11:12'; $ws.Cells.Item(164,3).Value = '15_ABASTO'; $ws.Cells.Item(164,4).Value = 46; $ws.Cells.Item(164,5).Value = 'LP1912'
$ws.Cells.Item(165,1).Value = '10:26:41'; $ws.Cells.Item(165,2).Value = '11:19'; $ws.Cells.Item(165,3).Value = '86_EST CHICA-ESC AGRARIA'; $ws.Cells.Item(165,4).Value = 53; $ws.Cells.Item(165,5).Value = 'LP1912'
$ws.Cells.Item(166,1).Value = '09:38:09'; $ws.Cells.Item(166,2).Value = '11:20'; $ws.Cells.Item(166,3).Value = '86_EST CHICA-ESC AGRARIA'; $ws.Cells.Item(166,4).Value = 102; $ws.Cells.Item(166,5).Value = 'LP1912'
$ws.Cells.Item(167,1).Value = '11:20:07'; $ws.Cells.Item(167,2).Value = '11:21'; $ws.Cells.Item(167,3).Value = '16_SANTA ANA'; $ws.Cells.Item(167,4).Value = 1; $ws.Cells.Item(167,5).Value = 'LP1912'
$ws.Cells.Item(168,1).Value = '09:38:09'; $ws.Cells.Item(168,2).Value = '11:21'; $ws.Cells.Item(168,3).Value = '26_HERNANDEZ'; $ws.Cells.Item(168,4).Value = 103; $ws.Cells.Item(168,5).Value = 'LP1912'
$ws.Cells.Item(169,1).Value = '11:20:07'; $ws.Cells.Item(169,2).Value = '11:22'; $ws.Cells.Item(169,3).Value = '17_ROMERO'; $ws.Cells.Item(169,4).Value = 2; $ws.Cells.Item(169,5).Value = 'LP1912'
$ws.Cells.Item(170,1).Value = '10:56:30'; $ws.Cells.Item(170,2).Value = '11:24'; $ws.Cells.Item(170,3).Value = '10_OLMOS'; $ws.Cells.Item(170,4).Value = 28; $ws.Cells.Item(170,5).Value = 'LP1912'
$ws.Cells.Item(171,1).Value = '11:20:07'; $ws.Cells.Item(171,2).Value = '11:25'; $ws.Cells.Item(171,3).Value = '16_SANTA ANA'; $ws.Cells.Item(171,4).Value = 5; $ws.Cells.Item(171,5).Value = 'LP1912'
$ws.Cells.Item(172,1).Value = '09:38:09'; $ws.Cells.Item(172,2).Value = '11:27'; $ws.Cells.Item(172,3).Value = '225_C ROCA-H SUR'; $ws.Cells.Item(172,4).Value = 109; $ws.Cells.Item(172,5).Value = 'LP1912'
$ws.Cells.Item(173,1).Value = '09:38:09'; $ws.Cells.Item(173,2).Value = '11:32'; $ws.Cells.Item(173,3).Value = '81_EL PELIGRO'; $ws.Cells.Item(173,4).Value = 114; $ws.Cells.Item(173,5).Value = 'LP1912'
$ws.Cells.Item(174,1).Value = '10:56:30'; $ws.Cells.Item(174,2).Value = '11:34'; $ws.Cells.Item(174,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(174,4).Value = 38; $ws.Cells.Item(174,5).Value = 'LP1912'
$ws.Cells.Item(175,1).Value = '11:20:07'; $ws.Cells.Item(175,2).Value = '11:35'; $ws.Cells.Item(175,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(175,4).Value = 15; $ws.Cells.Item(175,5).Value = 'LP1912'
$ws.Cells.Item(176,1).Value = '10:26:41'; $ws.Cells.Item(176,2).Value = '11:35'; $ws.Cells.Item(176,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(176,4).Value = 69; $ws.Cells.Item(176,5).Value = 'LP1912'
$ws.Cells.Item(177,1).Value = '09:38:09'; $ws.Cells.Item(177,2).Value = '11:36'; $ws.Cells.Item(177,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(177,4).Value = 118; $ws.Cells.Item(177,5).Value = 'LP1912'
$ws.Cells.Item(178,1).Value = '10:26:41'; $ws.Cells.Item(178,2).Value = '11:41'; $ws.Cells.Item(178,3).Value = '17_ROMERO'; $ws.Cells.Item(178,4).Value = 75; $ws.Cells.Item(178,5).Value = 'LP1912'
$ws.Cells.Item(179,1).Value = '10:56:30'; $ws.Cells.Item(179,2).Value = '11:42'; $ws.Cells.Item(179,3).Value = '17_ROMERO'; $ws.Cells.Item(179,4).Value = 46; $ws.Cells.Item(179,5).Value = 'LP1912'
$ws.Cells.Item(180,1).Value = '11:48:04'; $ws.Cells.Item(180,2).Value = '11:49'; $ws.Cells.Item(180,3).Value = '16_SANTA ANA'; $ws.Cells.Item(180,4).Value = 1; $ws.Cells.Item(180,5).Value = 'LP1912'
$ws.Cells.Item(181,1).Value = '10:26:41'; $ws.Cells.Item(181,2).Value = '11:51'; $ws.Cells.Item(181,3).Value = '215B_EL PATO'; $ws.Cells.Item(181,4).Value = 85; $ws.Cells.Item(181,5).Value = 'LP1912'
$ws.Cells.Item(182,1).Value = '10:56:30'; $ws.Cells.Item(182,2).Value = '11:52'; $ws.Cells.Item(182,3).Value = '15_ABASTO'; $ws.Cells.Item(182,4).Value = 56; $ws.Cells.Item(182,5).Value = 'LP1912'
$ws.Cells.Item(183,1).Value = '11:48:04'; $ws.Cells.Item(183,2).Value = '11:53'; $ws.Cells.Item(183,3).Value = '16_SANTA ANA'; $ws.Cells.Item(183,4).Value = 5; $ws.Cells.Item(183,5).Value = 'LP1912'
$ws.Cells.Item(184,1).Value = '10:26:41'; $ws.Cells.Item(184,2).Value = '11:59'; $ws.Cells.Item(184,3).Value = '225_GOMEZ'; $ws.Cells.Item(184,4).Value = 93; $ws.Cells.Item(184,5).Value = 'LP1912'
$ws.Cells.Item(185,1).Value = '12:01:50'; $ws.Cells.Item(185,2).Value = '12:01'; $ws.Cells.Item(185,3).Value = '16_SANTA ANA'; $ws.Cells.Item(185,4).Value = 0; $ws.Cells.Item(185,5).Value = 'LP1912'
$ws.Cells.Item(186,1).Value = '10:26:41'; $ws.Cells.Item(186,2).Value = '12:02'; $ws.Cells.Item(186,3).Value = '84_COLONIA URQUIZA-ESC 49'; $ws.Cells.Item(186,4).Value = 96; $ws.Cells.Item(186,5).Value = 'LP1912'
$ws.Cells.Item(187,1).Value = '12:01:50'; $ws.Cells.Item(187,2).Value = '12:04'; $ws.Cells.Item(187,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(187,4).Value = 3; $ws.Cells.Item(187,5).Value = 'LP1912'
$ws.Cells.Item(188,1).Value = '11:20:07'; $ws.Cells.Item(188,2).Value = '12:05'; $ws.Cells.Item(188,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(188,4).Value = 45; $ws.Cells.Item(188,5).Value = 'LP1912'
$ws.Cells.Item(189,1).Value = '10:26:41'; $ws.Cells.Item(189,2).Value = '12:06'; $ws.Cells.Item(189,3).Value = '16_P MOR-SANTA ANA'; $ws.Cells.Item(189,4).Value = 100; $ws.Cells.Item(189,5).Value = 'LP1912'
$ws.Cells.Item(190,1).Value = '10:56:30'; $ws.Cells.Item(190,2).Value = '12:06'; $ws.Cells.Item(190,3).Value = '14_ABASTO'; $ws.Cells.Item(190,4).Value = 70; $ws.Cells.Item(190,5).Value = 'LP1912'
$ws.Cells.Item(191,1).Value = '12:01:50'; $ws.Cells.Item(191,2).Value = '12:06'; $ws.Cells.Item(191,3).Value = '16_SANTA ANA'; $ws.Cells.Item(191,4).Value = 5; $ws.Cells.Item(191,5).Value = 'LP1912'
$ws.Cells.Item(192,1).Value = '12:01:50'; $ws.Cells.Item(192,2).Value = '12:06'; $ws.Cells.Item(192,3).Value = '84_COLONIA URQUIZA-ESC 49'; $ws.Cells.Item(192,4).Value = 5; $ws.Cells.Item(192,5).Value = 'LP1912'
$ws.Cells.Item(193,1).Value = '11:20:07'; $ws.Cells.Item(193,2).Value = '12:07'; $ws.Cells.Item(193,3).Value = '14_ABASTO'; $ws.Cells.Item(193,4).Value = 47; $ws.Cells.Item(193,5).Value = 'LP1912'
$ws.Cells.Item(194,1).Value = '11:20:07'; $ws.Cells.Item(194,2).Value = '12:07'; $ws.Cells.Item(194,3).Value = '16_P MOR-SANTA ANA'; $ws.Cells.Item(194,4).Value = 47; $ws.Cells.Item(194,5).Value = 'LP1912'
$ws.Cells.Item(195,1).Value = '10:56:30'; $ws.Cells.Item(195,2).Value = '12:10'; $ws.Cells.Item(195,3).Value = '10_OLMOS'; $ws.Cells.Item(195,4).Value = 74; $ws.Cells.Item(195,5).Value = 'LP1912'
$ws.Cells.Item(196,1).Value = '11:20:07'; $ws.Cells.Item(196,2).Value = '12:13'; $ws.Cells.Item(196,3).Value = '10_OLMOS'; $ws.Cells.Item(196,4).Value = 53; $ws.Cells.Item(196,5).Value = 'LP1912'
$ws.Cells.Item(197,1).Value = '10:26:41'; $ws.Cells.Item(197,2).Value = '12:14'; $ws.Cells.Item(197,3).Value = '17_ROMERO'; $ws.Cells.Item(197,4).Value = 108; $ws.Cells.Item(197,5).Value = 'LP1912'
$ws.Cells.Item(198,1).Value = '10:26:41'; $ws.Cells.Item(198,2).Value = '12:19'; $ws.Cells.Item(198,3).Value = '14_ABASTO'; $ws.Cells.Item(198,4).Value = 113; $ws.Cells.Item(198,5).Value = 'LP1912'
$ws.Cells.Item(199,1).Value = '10:56:30'; $ws.Cells.Item(199,2).Value = '12:20'; $ws.Cells.Item(199,3).Value = '14_ABASTO'; $ws.Cells.Item(199,4).Value = 84; $ws.Cells.Item(199,5).Value = 'LP1912'
$ws.Cells.Item(200,1).Value = '10:26:41'; $ws.Cells.Item(200,2).Value = '12:20'; $ws.Cells.Item(200,3).Value = '215A_EL PATO'; $ws.Cells.Item(200,4).Value = 114; $ws.Cells.Item(200,5).Value = 'LP1912'
$ws.Cells.Item(201,1).Value = '10:26:41'; $ws.Cells.Item(201,2).Value = '12:21'; $ws.Cells.Item(201,3).Value = '26_HERNANDEZ'; $ws.Cells.Item(201,4).Value = 115; $ws.Cells.Item(201,5).Value = 'LP1912'
$ws.Cells.Item(202,1).Value = '11:20:07'; $ws.Cells.Item(202,2).Value = '12:21'; $ws.Cells.Item(202,3).Value = '14_ABASTO'; $ws.Cells.Item(202,4).Value = 61; $ws.Cells.Item(202,5).Value = 'LP1912'
$ws.Cells.Item(203,1).Value = '11:20:07'; $ws.Cells.Item(203,2).Value = '12:21'; $ws.Cells.Item(203,3).Value = '215A_EL PATO'; $ws.Cells.Item(203,4).Value = 61; $ws.Cells.Item(203,5).Value = 'LP1912'
$ws.Cells.Item(204,1).Value = '12:01:50'; $ws.Cells.Item(204,2).Value = '12:34'; $ws.Cells.Item(204,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(204,4).Value = 33; $ws.Cells.Item(204,5).Value = 'LP1912'
$ws.Cells.Item(205,1).Value = '12:01:50'; $ws.Cells.Item(205,2).Value = '12:34'; $ws.Cells.Item(205,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(205,4).Value = 33; $ws.Cells.Item(205,5).Value = 'LP1912'
$ws.Cells.Item(206,1).Value = '11:48:04'; $ws.Cells.Item(206,2).Value = '12:35'; $ws.Cells.Item(206,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(206,4).Value = 47; $ws.Cells.Item(206,5).Value = 'LP1912'
$ws.Cells.Item(207,1).Value = '11:48:04'; $ws.Cells.Item(207,2).Value = '12:35'; $ws.Cells.Item(207,3).Value = '23_HERNANDEZ'; $ws.Cells.Item(207,4).Value = 47; $ws.Cells.Item(207,5).Value = 'LP1912'
$ws.Cells.Item(208,1).Value = '10:56:30'; $ws.Cells.Item(208,2).Value = '12:36'; $ws.Cells.Item(208,3).Value = '27_EL RETIRO'; $ws.Cells.Item(208,4).Value = 100; $ws.Cells.Item(208,5).Value = 'LP1912'
$ws.Cells.Item(209,1).Value = '11:20:07'; $ws.Cells.Item(209,2).Value = '12:37'; $ws.Cells.Item(209,3).Value = '27_EL RETIRO'; $ws.Cells.Item(209,4).Value = 77; $ws.Cells.Item(209,5).Value = 'LP1912'
$ws.Cells.Item(210,1).Value = '10:56:30'; $ws.Cells.Item(210,2).Value = '12:38'; $ws.Cells.Item(210,3).Value = '17_179 Y 38'; $ws.Cells.Item(210,4).Value = 102; $ws.Cells.Item(210,5).Value = 'LP1912'
$ws.Cells.Item(211,1).Value = '10:56:30'; $ws.Cells.Item(211,2).Value = '12:41'; $ws.Cells.Item(211,3).Value = '10_OLMOS'; $ws.Cells.Item(211,4).Value = 105; $ws.Cells.Item(211,5).Value = 'LP1912'
$ws.Cells.Item(212,1).Value = '12:01:50'; $ws.Cells.Item(212,2).Value = '12:48'; $ws.Cells.Item(212,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(212,4).Value = 47; $ws.Cells.Item(212,5).Value = 'LP1912'
$ws.Cells.Item(213,1).Value = '11:20:07'; $ws.Cells.Item(213,2).Value = '12:49'; $ws.Cells.Item(213,3).Value = '11_ETCHEVERRY'; $ws.Cells.Item(213,4).Value = 89; $ws.Cells.Item(213,5).Value = 'LP1912'
$ws.Cells.Item(214,1).Value = '12:01:50'; $ws.Cells.Item(214,2).Value = '12:56'; $ws.Cells.Item(214,3).Value = '14_ABASTO'; $ws.Cells.Item(214,4).Value = 55; $ws.Cells.Item(214,5).Value = 'LP1912'
$ws.Cells.Item(215,1).Value = '11:20:07'; $ws.Cells.Item(215,2).Value = '13:02'; $ws.Cells.Item(215,3).Value = '15_ABASTO'; $ws.Cells.Item(215,4).Value = 102; $ws.Cells.Item(215,5).Value = 'LP1912'
$ws.Cells.Item(216,1).Value = '11:20:07'; $ws.Cells.Item(216,2).Value = '13:07'; $ws.Cells.Item(216,3).Value = '16_P MOR-SANTA ANA'; $ws.Cells.Item(216,4).Value = 107; $ws.Cells.Item(216,5).Value = 'LP1912'
$ws.Cells.Item(217,1).Value = '12:01:50'; $ws.Cells.Item(217,2).Value = '13:13'; $ws.Cells.Item(217,3).Value = '215D_EL PATO'; $ws.Cells.Item(217,4).Value = 72; $ws.Cells.Item(217,5).Value = 'LP1912'
$ws.Cells.Item(218,1).Value = '11:20:07'; $ws.Cells.Item(218,2).Value = '13:14'; $ws.Cells.Item(218,3).Value = '215D_EL PATO'; $ws.Cells.Item(218,4).Value = 114; $ws.Cells.Item(218,5).Value = 'LP1912'
$ws.Cells.Item(219,1).Value = '12:01:50'; $ws.Cells.Item(219,2).Value = '13:19'; $ws.Cells.Item(219,3).Value = '10_OLMOS'; $ws.Cells.Item(219,4).Value = 78; $ws.Cells.Item(219,5).Value = 'LP1912'
$ws.Cells.Item(220,1).Value = '11:48:04'; $ws.Cells.Item(220,2).Value = '13:20'; $ws.Cells.Item(220,3).Value = '10_OLMOS'; $ws.Cells.Item(220,4).Value = 92; $ws.Cells.Item(220,5).Value = 'LP1912'
$ws.Cells.Item(221,1).Value = '11:48:04'; $ws.Cells.Item(221,2).Value = '13:21'; $ws.Cells.Item(221,3).Value = '26_HERNANDEZ'; $ws.Cells.Item(221,4).Value = 93; $ws.Cells.Item(221,5).Value = 'LP1912'
$ws.Cells.Item(222,1).Value = '11:48:04'; $ws.Cells.Item(222,2).Value = '13:27'; $ws.Cells.Item(222,3).Value = '14_ABASTO'; $ws.Cells.Item(222,4).Value = 99; $ws.Cells.Item(222,5).Value = 'LP1912'
$ws.Cells.Item(223,1).Value = '11:48:04'; $ws.Cells.Item(223,2).Value = '13:36'; $ws.Cells.Item(223,3).Value = '15_ABASTO'; $ws.Cells.Item(223,4).Value = 108; $ws.Cells.Item(223,5).Value = 'LP1912'
$ws.Cells.Item(224,1).Value = '11:48:04'; $ws.Cells.Item(224,2).Value = '13:46'; $ws.Cells.Item(224,3).Value = '17_ROMERO'; $ws.Cells.Item(224,4).Value = 118; $ws.Cells.Item(224,5).Value = 'LP1912'
$ws.Cells.Item(225,1).Value = '12:01:50'; $ws.Cells.Item(225,2).Value = '13:50'; $ws.Cells.Item(225,3).Value = '215A_EL PATO'; $ws.Cells.Item(225,4).Value = 109; $ws.Cells.Item(225,5).Value = 'LP1912'
$ws.Cells.Item(226,1).Value = '12:01:50'; $ws.Cells.Item(226,2).Value = '13:52'; $ws.Cells.Item(226,3).Value = '10_OLMOS'; $ws.Cells.Item(226,4).Value = 111; $ws.Cells.Item(226,5).Value = 'LP1912'
$ws.Cells.Item(227,1).Value = '12:01:50'; $ws.Cells.Item(227,2).Value = '13:55'; $ws.Cells.Item(227,3).Value = '225_GOMEZ'; $ws.Cells.Item(227,4).Value = 114; $ws.Cells.Item(227,5).Value = 'LP1912'
$ws.Cells.Item(228,1).Value = '12:01:50'; $ws.Cells.Item(228,2).Value = '13:56'; $ws.Cells.Item(228,3).Value = '16_P MOR-167 Y 521'; $ws.Cells.Item(228,4).Value = 115; $ws.Cells.Item(228,5).Value = 'LP1912'

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Cells.Item(2,1).Value = 'Última actualización: 12:01:50'
$ws.Cells.Item(3,1).Value = 'Total filas: 27'

$ws.Cells.Item(30,1).Value = '12:01:50'; $ws.Cells.Item(30,2).Value = '13:13'; $ws.Cells.Item(30,3).Value = '215D_EL PATO'; $ws.Cells.Item(30,4).Value = 72; $ws.Cells.Item(30,5).Value = 'LP1912'
$ws.Cells.Item(31,1).Value = '11:20:07'; $ws.Cells.Item(31,2).Value = '13:14'; $ws.Cells.Item(31,3).Value = '215D_EL PATO'; $ws.Cells.Item(31,4).Value = 114; $ws.Cells.Item(31,5).Value = 'LP1912'
$ws.Cells.Item(32,1).Value = '12:01:50'; $ws.Cells.Item(32,2).Value = '13:50'; $ws.Cells.Item(32,3).Value = '215A_EL PATO'; $ws.Cells.Item(32,4).Value = 109; $ws.Cells.Item(32,5).Value = 'LP1912'

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Cells.Item(2,1).Value = 'Última actualización: 12:01:50'
$ws.Cells.Item(3,1).Value = 'Total filas: 37'

$ws.Cells.Item(37,1).Value = '12:01:50'; $ws.Cells.Item(37,2).Value = '12:05'; $ws.Cells.Item(37,3).Value = '215A_LA PLATA'; $ws.Cells.Item(37,4).Value = 4; $ws.Cells.Item(37,5).Value = 'L6173'
$ws.Cells.Item(38,1).Value = '11:48:04'; $ws.Cells.Item(38,2).Value = '12:06'; $ws.Cells.Item(38,3).Value = '215A_LA PLATA'; $ws.Cells.Item(38,4).Value = 18; $ws.Cells.Item(38,5).Value = 'L6173'
$ws.Cells.Item(39,1).Value = '12:01:50'; $ws.Cells.Item(39,2).Value = '12:53'; $ws.Cells.Item(39,3).Value = '215C_LA PLATA'; $ws.Cells.Item(39,4).Value = 52; $ws.Cells.Item(39,5).Value = 'L6203'
$ws.Cells.Item(40,1).Value = '10:56:30'; $ws.Cells.Item(40,2).Value = '12:54'; $ws.Cells.Item(40,3).Value = '215C_LA PLATA'; $ws.Cells.Item(40,4).Value = 118; $ws.Cells.Item(40,5).Value = 'L6203'
$ws.Cells.Item(41,1).Value = '12:01:50'; $ws.Cells.Item(41,2).Value = '13:30'; $ws.Cells.Item(41,3).Value = '215B_LP-P MOR-1 Y 57'; $ws.Cells.Item(41,4).Value = 89; $ws.Cells.Item(41,5).Value = 'L6173'
$ws.Cells.Item(42,1).Value = '11:48:04'; $ws.Cells.Item(42,2).Value = '13:31'; $ws.Cells.Item(42,3).Value = '215B_LP-P MOR-1 Y 57'; $ws.Cells.Item(42,4).Value = 103; $ws.Cells.Item(42,5).Value = 'LP1912'

